# ---------------------------------------------------------------------------
# SGX 설치.docx edit script
#
# 1. Remove the stray "_GoBack" bookmark pair that originally sat inside the
#    right-aligned tab paragraph near the top of the document.
# 2. Append a short "command line" walkthrough (one intro sentence plus six
#    systemctl commands) right after the screenshot that illustrates step
#    25 ("자동 소프트웨어 업데이트 중단 시키기"), i.e. in the gap between
#    that picture and the following "26. 예제 코드 실행하기" heading. The
#    "_GoBack" bookmark re-appears at the end of the new block (Word always
#    keeps exactly one), followed by one trailing empty paragraph.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Step 1: drop the old _GoBack bookmark -------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
if ($goBack -ne $null) {
    [void]$goBack.Delete()
}

# --- Step 2: locate the empty paragraph right after the step-25 picture --
$paragraphs = $d.Paragraphs
$insertAfterIndex = -1
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $paragraphs.Item($i)
    if ($candidate.Range.InlineShapes.Count -gt 0) {
        $insertAfterIndex = $i + 1
    }
}

$targetParagraph = $paragraphs.Item($insertAfterIndex)

# Create one fresh empty paragraph right after it, then expand that single
# paragraph into the full new block via InsertXML (keeps the original
# empty placeholder paragraph untouched).
[void]$targetParagraph.Range.InsertParagraphAfter()
$newParagraph = $d.Paragraphs.Item($insertAfterIndex + 1)

$newBlockXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>커맨드라인으로 수행할 경우에는 아래를 타이핑한다.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>systemctl stop apt-daily.timer</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>systemctl disable apt-daily.timer</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>systemctl disable apt-daily.service</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>systemctl stop apt-daily-upgrade.timer</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>systemctl disable apt-daily-upgrade.timer</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>systemctl disable apt-daily-upgrade.service</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$newParagraph.Range.InsertXML($newBlockXml)
